$d = $word.ActiveDocument
$d.TrackRevisions = $false

# 1. HPCDATAMGM-1721 -> HPCDATAMGM-1745
$d.Content.Find.Execute("HPCDATAMGM-1721", $true, $false, $false, $false, $false, $true, 1, $false, "HPCDATAMGM-1745", 2)

# 2. "the owners or co-owners of the file" -> "the owner or co-owner of the file"
$d.Content.Find.Execute("the owners or co-owners of the file", $true, $false, $false, $false, $false, $true, 1, $false, "the owner or co-owner of the file", 2)

# 3. CLU command description wording change
$d.Content.Find.Execute("to obtain the presigned URL to download a file to.", $true, $false, $false, $false, $false, $true, 1, $false, "to obtain the presigned download URL for a file.", 2)

# 4. "Google drive upload" -> "Google Drive upload"
$d.Content.Find.Execute("Google drive upload", $true, $false, $false, $false, $false, $true, 1, $false, "Google Drive upload", 2)
